$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = "  -0.26%  "
    3 = "  +1.32%  "
    4 = "  -0.10%  "
    5 = "  +1.02%  "
    6 = "  +0.16%  "
    7 = "  -0.05%  "
    8 = "  -0.20%  "
    9 = "  -4.92%  "
    10 = "  -1.03%  "
    11 = "  +3.01%  "
    12 = "  +1.37%  "
    13 = "  -3.23%  "
    15 = "  +2.66%  "
    16 = "  -0.05%  "
    17 = "  -1.23%  "
    18 = "  +1.29%  "
    19 = "  +4.39%  "
    20 = "  +2.34%  "
    21 = "  +0.40%  "
    22 = "  +0.02%  "
    23 = "  +1.46%  "
    24 = "  +0.64%  "
    25 = "  +1.40%  "
    26 = "  -0.08%  "
    27 = "  -2.62%  "
    28 = "  -1.69%  "
    29 = "  -1.16%  "
    30 = "  -1.56%  "
    31 = "  -0.35%  "
    32 = "  +0.56%  "
    33 = "  +0.85%  "
    34 = "  -1.83%  "
    35 = "  -0.61%  "
    36 = "  +0.25%  "
    37 = "  -0.73%  "
    38 = "  +1.79%  "
    39 = "  +0.38%  "
    40 = "  -1.21%  "
    41 = "  +0.89%  "
    42 = "  +2.61%  "
    43 = "  -0.06%  "
    44 = "  -1.05%  "
    45 = "  -1.13%  "
    46 = "  +0.90%  "
    47 = "  -3.67%  "
    48 = "  +1.58%  "
    49 = "  +3.43%  "
    50 = "  +7.50%  "
    51 = "  -0.30%  "
}

foreach ($row in $updates.Keys) {
    $ws.Range("E$row").Value = $updates[$row]
}

Write-Output "Updated $($updates.Count) cells in column E"
